$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style (default, unstyled data cell) used to force text-typed
# values for cells whose new content would otherwise auto-parse as a number,
# so the saved cell keeps the original inlineStr/text semantics.
$plainStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = "64.745.10"
$ws.Range("D3").Value = "3.462.08"
$ws.Range("E3").Value = "  +3.91%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.09"
$ws.Range("D5").Style = $plainStyle
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.55"
$ws.Range("D6").Style = $plainStyle
$ws.Range("E6").Value = "  +3.00%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "3.459.03"
$ws.Range("E8").Value = "  +3.43%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.556"
$ws.Range("D9").Style = $plainStyle
$ws.Range("E9").Value = "  +5.34%  "
$ws.Range("E10").Value = "  +0.98%  "
$ws.Range("E11").Value = "  +6.32%  "
$ws.Range("E12").Value = "  +1.98%  "
$ws.Range("D13").Value = "4.057.06"
$ws.Range("E13").Value = "  +3.95%  "
$ws.Range("E14").Value = "  -1.57%  "
$ws.Range("E15").Value = "  +9.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.91"
$ws.Range("D16").Style = $plainStyle
$ws.Range("E16").Value = "  +3.14%  "
$ws.Range("D17").Value = "64.671.62"
$ws.Range("E17").Value = "  +3.57%  "
$ws.Range("D18").Value = "3.487.28"
$ws.Range("E18").Value = "  +4.77%  "
$ws.Range("E19").Value = "  -0.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.39"
$ws.Range("D20").Style = $plainStyle
$ws.Range("E20").Value = "  +4.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "398.55"
$ws.Range("D21").Style = $plainStyle
$ws.Range("E21").Value = "  +3.63%  "
$ws.Range("E22").Value = "  +0.75%  "
$ws.Range("E24").Value = "  +3.11%  "
$ws.Range("E25").Value = "  -0.87%  "
$ws.Range("E26").Value = "  +23.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.51"
$ws.Range("D27").Style = $plainStyle
$ws.Range("E27").Value = "  +6.29%  "
$ws.Range("E28").Value = "  +3.14%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.99"
$ws.Range("D30").Style = $plainStyle
$ws.Range("E30").Value = "  +8.32%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.05"
$ws.Range("D31").Style = $plainStyle
$ws.Range("E31").Value = "  +4.24%  "
$ws.Range("B32").Value = "RenderToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.70"
$ws.Range("D32").Style = $plainStyle
$ws.Range("E32").Value = "  +5.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.39"
$ws.Range("D33").Style = $plainStyle
$ws.Range("E33").Value = "  +5.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.87"
$ws.Range("D34").Style = $plainStyle
$ws.Range("E34").Value = "  +3.88%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.01"
$ws.Range("D36").Style = $plainStyle
$ws.Range("E36").Value = "  +3.46%  "
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "161.39"
$ws.Range("D38").Style = $plainStyle
$ws.Range("E38").Value = "  +1.26%  "
$ws.Range("E39").Value = "  +7.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "27.75"
$ws.Range("D40").Style = $plainStyle
$ws.Range("E40").Value = "  +2.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.88"
$ws.Range("D41").Style = $plainStyle
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").Value = "2.901.14"
$ws.Range("E42").Value = "  +1.43%  "
$ws.Range("E43").Value = "  +2.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.776"
$ws.Range("D44").Style = $plainStyle
$ws.Range("E44").Value = "  +3.44%  "
$ws.Range("E45").Value = "  +1.94%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "41.89"
$ws.Range("D46").Style = $plainStyle
$ws.Range("E46").Value = "  +2.66%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.64"
$ws.Range("D47").Style = $plainStyle
$ws.Range("E47").Value = "  +7.31%  "
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.10"
$ws.Range("D48").Style = $plainStyle
$ws.Range("E48").Value = "  +4.86%  "
$ws.Range("E49").Value = "  +22.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.866"
$ws.Range("D50").Style = $plainStyle
$ws.Range("E50").Value = "  +6.48%  "
$ws.Range("E51").Value = "  +4.14%  "
